# "added search by name"
# The sheet held 3 alumni records (Arina Ixescu, Marin Marinescu, Maria
# Cutarescu). The edit replaces the single remaining record (row 2) with a
# new person (Alexandru Cutarescu) and removes the other two records
# (rows 3 and 4) entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 in place with the new record ------------------------
$ws.Range("A2").Value = "3215515848587"
$ws.Range("B2").Value = "Alexandru"
$ws.Range("C2").Value = "Cutarescu"
$ws.Range("D2").Value = "alexandru.cutarescu@gmail.com"
$ws.Range("E2").Value = 2020
$ws.Range("F2").Value = "Angular"

# --- Drop the two rows that are no longer needed -----------------------
$ws.Rows("3:4").Delete()

# --- Hyperlinks: the old mailto: links (for the now-deleted rows 3 & 4,
# plus the stale one on row 2) need to be replaced with a single mailto:
# link that points at the new e-mail address in D2.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:alexandru.cutarescu@gmail.com") | Out-Null
$ws.Range("D2").Style = "Гиперссылка"

# --- Leave the selection where the author left it in the saved file ----
$ws.Range("D3").Select()
